$d = $word.ActiveDocument

$replacements = @(
    @{old = "406÷4="; new = "423÷8="},
    @{old = "410÷6="; new = "430÷5="},
    @{old = "978÷8="; new = "526÷6="},
    @{old = "848÷3="; new = "809÷6="},
    @{old = "815÷8="; new = "162÷6="},
    @{old = "331÷7="; new = "269÷6="},
    @{old = "997÷3="; new = "728÷4="},
    @{old = "195÷2="; new = "626÷4="},
    @{old = "962÷5="; new = "944÷4="},
    @{old = "508÷7="; new = "817÷2="},
    @{old = "480÷7="; new = "925÷6="},
    @{old = "831÷6="; new = "196÷4="},
    @{old = "522÷3="; new = "381÷8="},
    @{old = "272÷2="; new = "528÷3="},
    @{old = "906÷2="; new = "336÷8="},
    @{old = "648÷3="; new = "429÷6="},
    @{old = "470÷5="; new = "941÷7="},
    @{old = "872÷8="; new = "383÷7="},
    @{old = "989÷6="; new = "384÷4="},
    @{old = "828÷5="; new = "199÷5="},
    @{old = "543÷5="; new = "763÷4="},
    @{old = "178÷6="; new = "218÷4="},
    @{old = "438÷6="; new = "165÷6="},
    @{old = "725÷3="; new = "367÷2="},
    @{old = "587÷8="; new = "316÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
